$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.249.11"
$ws.Range("E2").Value = "  +10.71%  "

# Row 3
$ws.Range("D3").Value = "3.472.47"
$ws.Range("E3").Value = "  +6.61%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.66"
$ws.Range("E5").Value = "  +4.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.58"
$ws.Range("E6").Value = "  +13.53%  "

# Row 7
$ws.Range("D7").Value = "3.464.56"
$ws.Range("E7").Value = "  +6.56%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +4.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.687"
$ws.Range("E10").Value = "  +10.68%  "

# Row 11
$ws.Range("E11").Value = "  +37.54%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.34"
$ws.Range("E12").Value = "  +5.31%  "

# Row 13
$ws.Range("E13").Value = "  -0.29%  "

# Row 14
$ws.Range("D14").Value = "4.026.44"
$ws.Range("E14").Value = "  +6.72%  "

# Row 15
$ws.Range("E15").Value = "  +4.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.98"
$ws.Range("E16").Value = "  +5.38%  "

# Row 17
$ws.Range("D17").Value = "3.469.18"
$ws.Range("E17").Value = "  +6.27%  "

# Row 18
$ws.Range("D18").Value = "63.065.24"
$ws.Range("E18").Value = "  +10.76%  "

# Row 19
$ws.Range("E19").Value = "  +0.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.02"
$ws.Range("E20").Value = "  -0.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000142"
$ws.Range("E21").Value = "  +31.60%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.32"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.66"
$ws.Range("E23").Value = "  +11.79%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "314.11"
$ws.Range("E24").Value = "  +6.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.89"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.18"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.02"
$ws.Range("E27").Value = "  +10.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  +4.50%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.87"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.175"
$ws.Range("E30").Value = "  +3.84%  "

# Row 31
$ws.Range("E31").Value = "  -1.56%  "

# Row 32
$ws.Range("E32").Value = "  +2.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.59"
$ws.Range("E33").Value = "  +3.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("E34").Value = "  +19.35%  "

# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.82"
$ws.Range("E36").Value = "  +2.68%  "

# Row 37
$ws.Range("E37").Value = "  -0.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.12"
$ws.Range("E38").Value = "  +1.65%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.46"
$ws.Range("E40").Value = "  +0.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +1.27%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.127"
$ws.Range("E42").Value = "  +4.63%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").Value = "  +7.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.22"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.14"
$ws.Range("E45").Value = "  +2.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.284"
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
$ws.Range("E48").Value = "  +1.83%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.03"
$ws.Range("E49").Value = "  -1.44%  "

# Row 50
$ws.Range("D50").Value = "3.819.67"

# Row 51
$ws.Range("D51").Value = "2.182.29"
$ws.Range("E51").Value = "  +1.80%  "
